$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 2550
$ws.Range("I18").Value = 2850
$ws.Range("K18").Value = 2850
$ws.Range("M18").Value = -2566

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1826.7778
$ws.Range("I45").Value = 1826.7778
$ws.Range("K45").Value = 1826.7778
$ws.Range("M45").Value = -1449.7778

# Row 122
$ws.Range("H122").Value = 1596.3334
$ws.Range("I122").Value = 1619.6
$ws.Range("K122").Value = 4858.799999999999
$ws.Range("M122").Value = -2408.799999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 3000
$ws.Range("I7").Value = 3000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2887
$ws.Range("N7").ClearContents()

# Row 134
$ws.Range("H134").Value = 3583.7026
$ws.Range("I134").Value = 909.0303
$ws.Range("J134").Value = 25649.75
$ws.Range("K134").Value = 2727.0909
$ws.Range("L134").Value = 76949.25
$ws.Range("M134").Value = -192.0909000000001
$ws.Range("N134").Value = -82019.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1175.6111
$ws.Range("I31").Value = 1209.4706
$ws.Range("J31").Value = 600
$ws.Range("K31").Value = 1209.4706
$ws.Range("L31").Value = 600
$ws.Range("M31").Value = -914.4706000000001
$ws.Range("N31").Value = -1190

# Row 34
$ws.Range("H34").Value = 1175.6111
$ws.Range("I34").Value = 1209.4706
$ws.Range("J34").Value = 600
$ws.Range("K34").Value = 1209.4706
$ws.Range("L34").Value = 600
$ws.Range("M34").Value = -1007.4706
$ws.Range("N34").Value = -1004

# Row 58
$ws.Range("H58").Value = 660.875
$ws.Range("I58").Value = 648.61536
$ws.Range("J58").Value = 714
$ws.Range("K58").Value = 648.61536
$ws.Range("L58").Value = 714
$ws.Range("M58").Value = -445.61536
$ws.Range("N58").Value = -1120

# Row 122
$ws.Range("H122").Value = 1033.2
$ws.Range("I122").Value = 991.5
$ws.Range("K122").Value = 2974.5
$ws.Range("M122").Value = -524.5

# Row 132
$ws.Range("H132").Value = 7636.8423
$ws.Range("I132").Value = 9704.416999999999
$ws.Range("K132").Value = 29113.251
$ws.Range("M132").Value = -26583.251

# Row 136
$ws.Range("H136").Value = 660.875
$ws.Range("I136").Value = 648.61536
$ws.Range("J136").Value = 714
$ws.Range("K136").Value = 1945.84608
$ws.Range("L136").Value = 2142
$ws.Range("M136").Value = 604.15392
$ws.Range("N136").Value = -7242

$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 1444.5
$ws.Range("I86").Value = 1444.5
$ws.Range("K86").Value = 4333.5
$ws.Range("M86").Value = -3147.5

# Row 89
$ws.Range("H89").Value = 1444.5
$ws.Range("I89").Value = 1444.5
$ws.Range("K89").Value = 13000.5
$ws.Range("M89").Value = -7072.5

# Row 131
$ws.Range("H131").Value = 20001372
$ws.Range("J131").Value = 1591.2195
$ws.Range("L131").Value = 4773.6585
$ws.Range("N131").Value = -14853.6585

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 16671471
$ws.Range("I70").Value = 14710095
$ws.Range("K70").Value = 14710095
$ws.Range("M70").Value = -14709825

# Row 73
$ws.Range("H73").Value = 16671471
$ws.Range("I73").Value = 14710095
$ws.Range("K73").Value = 14710095
$ws.Range("M73").Value = -14709159

# Row 122
$ws.Range("H122").Value = 1162.25
$ws.Range("I122").Value = 1049.6666
$ws.Range("K122").Value = 3148.9998
$ws.Range("M122").Value = -698.9998000000001

# Row 123
$ws.Range("H123").Value = 10354.214
$ws.Range("J123").Value = 10354.214
$ws.Range("L123").Value = 10354.214
$ws.Range("N123").Value = -15254.214

# Row 128
$ws.Range("H128").Value = 38790.363
$ws.Range("J128").Value = 38790.363
$ws.Range("L128").Value = 38790.363
$ws.Range("N128").Value = -48750.363

# Row 132
$ws.Range("H132").Value = 2496.0625
$ws.Range("I132").Value = 2462.077
$ws.Range("K132").Value = 7386.231000000001
$ws.Range("M132").Value = -4856.231000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1505.4546
$ws.Range("I22").Value = 1397.7778
$ws.Range("J22").Value = 1990
$ws.Range("K22").Value = 1397.7778
$ws.Range("L22").Value = 1990
$ws.Range("M22").Value = -1102.7778
$ws.Range("N22").Value = -2580

# Row 27
$ws.Range("H27").Value = 1505.4546
$ws.Range("I27").Value = 1397.7778
$ws.Range("J27").Value = 1990
$ws.Range("K27").Value = 1397.7778
$ws.Range("L27").Value = 1990
$ws.Range("M27").Value = -1290.7778
$ws.Range("N27").Value = -2204

# Row 46
$ws.Range("H46").Value = 7858.3335
$ws.Range("I46").Value = 1100.3334
$ws.Range("J46").Value = 10111
$ws.Range("K46").Value = 1100.3334
$ws.Range("L46").Value = 10111
$ws.Range("M46").Value = -912.3334
$ws.Range("N46").Value = -10487

# Row 55
$ws.Range("H55").Value = 613.6667
$ws.Range("I55").Value = 336.4
$ws.Range("J55").Value = 2000
$ws.Range("K55").Value = 336.4
$ws.Range("L55").Value = 2000
$ws.Range("M55").Value = -163.4
$ws.Range("N55").Value = -2346

# Row 122
$ws.Range("H122").Value = 20836466
$ws.Range("I122").Value = 41669732
$ws.Range("J122").Value = 3201.3333
$ws.Range("K122").Value = 125009196
$ws.Range("L122").Value = 9603.999899999999
$ws.Range("M122").Value = -125006746
$ws.Range("N122").Value = -14503.9999

# Row 136
$ws.Range("H136").Value = 18385.334
$ws.Range("I136").Value = 26328
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 78984
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -76434
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 2340
$ws.Range("I96").Value = 2280
$ws.Range("J96").Value = 2370
$ws.Range("K96").Value = 2280
$ws.Range("L96").Value = 2370
$ws.Range("M96").Value = -907
$ws.Range("N96").Value = -5116

# Row 122
$ws.Range("H122").Value = 10404513
$ws.Range("I122").Value = 12385868
$ws.Range("J122").Value = 2399.75
$ws.Range("K122").Value = 37157604
$ws.Range("L122").Value = 7199.25
$ws.Range("M122").Value = -37155154
$ws.Range("N122").Value = -12099.25

# Row 136
$ws.Range("H136").Value = 473.16666
$ws.Range("I136").Value = 327.8
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 983.4000000000001
$ws.Range("L136").Value = 3600
$ws.Range("M136").Value = 1566.6
$ws.Range("N136").Value = -8700

Write-Host "Applied Kujata_Profits market price updates across all sheets."